$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row cells: "_old" -> "_FV2310", "_new" -> "_FV2404"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2310")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2404")
    }
}

# Freeze the header row (row 1)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the full used data range into a native Excel table
$range = $ws.UsedRange
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"
